# Weekly update: insert a new data row at the top of the date-ordered block
# (row 213), pushing the existing rows 213:225 down to 214:226.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 213 - everything below (213:225) shifts down
# to (214:226), matching the diff exactly.
$ws.Rows("213:213").Insert()

# The new row 213 repeats the same market/category/quality metadata as the
# surrounding rows, with a new date and new price figures.
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 44783
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 100112037
$ws.Cells.Item(213, 7).Value = "Cebollín"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 1600
$ws.Cells.Item(213, 11).Value = 1400
$ws.Cells.Item(213, 12).Value = 1600
$ws.Cells.Item(213, 13).Value = 1500
$ws.Cells.Item(213, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(213, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(213, 16).Value = 250
$ws.Cells.Item(213, 17).Value = 6
$ws.Cells.Item(213, 18).Value = "Hortaliza"
